$wb = $excel.ActiveWorkbook

# Map of cell -> new value, identical for both the "展览" and "全部类型" sheets
$updates = @{
    "F2"  = 9
    "F3"  = 5011
    "F5"  = 7273
    "F13" = 1709
    "F16" = 2853
    "F20" = 455
    "F23" = 275
    "F28" = 1335
    "F36" = 2671
    "F37" = 687
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
